$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D2:E51 retain Text format so numeric-looking strings are not
# auto-converted to numbers by Excel when the .Value is assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '56.880.72'
$ws.Range('E2').Value = '  +4.12%  '
$ws.Range('D3').Value = '2.343.71'
$ws.Range('E3').Value = '  +2.87%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').Value = '518.97'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').Value = '133.97'
$ws.Range('E6').Value = '  +4.26%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.535'
$ws.Range('E8').Value = '  +1.43%  '
$ws.Range('D9').Value = '2.342.38'
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('E10').Value = '  +7.23%  '
$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '5.26'
$ws.Range('E12').Value = '  +7.03%  '
$ws.Range('D13').Value = '0.340'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '23.73'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '2.740.94'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '56.764.33'
$ws.Range('E16').Value = '  +3.82%  '
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '2.330.95'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('D19').Value = '10.42'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('D20').Value = '4.25'
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').Value = '321.02'
$ws.Range('E21').Value = '  +4.61%  '
$ws.Range('D22').Value = '6.56'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = '60.45'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = '0.161'
$ws.Range('E26').Value = '  +7.08%  '
$ws.Range('E27').Value = '  +4.14%  '
$ws.Range('D28').Value = '1.23'
$ws.Range('E28').Value = '  +10.28%  '
$ws.Range('D29').Value = '170.18'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('E30').Value = '  +5.68%  '
$ws.Range('D31').Value = '1.69'
$ws.Range('E31').Value = '  +4.26%  '
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('D33').Value = '18.26'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').Value = '0.932'
$ws.Range('E36').Value = '  +2.22%  '
$ws.Range('E37').Value = '  +3.93%  '
$ws.Range('E38').Value = '  +5.30%  '
$ws.Range('D39').Value = '1.52'
$ws.Range('E39').Value = '  +7.71%  '
$ws.Range('D40').Value = '37.57'
$ws.Range('E40').Value = '  +2.95%  '
$ws.Range('D41').Value = '0.381'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('E42').Value = '  +5.94%  '
$ws.Range('D43').Value = '137.30'
$ws.Range('E43').Value = '  +8.51%  '
$ws.Range('D44').Value = '275.76'
$ws.Range('E44').Value = '  +10.17%  '
$ws.Range('D45').Value = '5.12'
$ws.Range('E45').Value = '  +6.18%  '
$ws.Range('D46').Value = '0.0931'
$ws.Range('E46').Value = '  +3.33%  '
$ws.Range('D47').Value = '0.0505'
$ws.Range('E47').Value = '  +1.94%  '
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('E49').Value = '  +5.06%  '
$ws.Range('D50').Value = '0.380'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '16.84'
$ws.Range('E51').Value = '  +2.63%  '
